$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '66.902.61'
Set-TextValue 'E2' '  -3.74%  '
Set-TextValue 'D3' '3.519.61'
Set-TextValue 'E3' '  -4.20%  '
Set-TextValue 'E4' '  +0.12%  '
Set-TextValue 'D5' '609.33'
Set-TextValue 'E5' '  -5.66%  '
Set-TextValue 'D6' '151.95'
Set-TextValue 'E6' '  -4.64%  '
Set-TextValue 'D7' '3.519.44'
Set-TextValue 'E7' '  -4.25%  '
Set-TextValue 'E8' '  +0.11%  '
Set-TextValue 'E9' '  -3.38%  '
Set-TextValue 'E10' '  -4.07%  '
Set-TextValue 'D11' '6.83'
Set-TextValue 'E11' '  -4.47%  '
Set-TextValue 'D12' '0.425'
Set-TextValue 'E12' '  -3.80%  '
Set-TextValue 'D13' '0.0000218'
Set-TextValue 'E13' '  -5.56%  '
Set-TextValue 'D14' '4.119.85'
Set-TextValue 'E14' '  -3.92%  '
Set-TextValue 'D15' '31.52'
Set-TextValue 'E15' '  -3.40%  '
Set-TextValue 'D16' '3.526.16'
Set-TextValue 'E16' '  -4.31%  '
Set-TextValue 'D17' '66.892.98'
Set-TextValue 'E17' '  -3.71%  '
Set-TextValue 'E18' '  +0.49%  '
Set-TextValue 'D19' '6.27'
Set-TextValue 'E19' '  -3.30%  '
Set-TextValue 'D20' '15.31'
Set-TextValue 'E20' '  -4.06%  '
Set-TextValue 'D21' '442.41'
Set-TextValue 'E21' '  -6.18%  '
Set-TextValue 'D22' '9.17'
Set-TextValue 'E22' '  -8.86%  '
Set-TextValue 'D23' '0.628'
Set-TextValue 'E23' '  -3.10%  '
Set-TextValue 'D24' '77.65'
Set-TextValue 'E24' '  -2.22%  '
Set-TextValue 'B25' 'Dai'
Set-TextValue 'C25' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D25' '1.00'
Set-TextValue 'E25' '  -0.08%  '
Set-TextValue 'B26' 'WrappedeETH'
Set-TextValue 'C26' 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue 'D26' '3.664.30'
Set-TextValue 'E26' '  -4.02%  '
Set-TextValue 'E27' '  -1.93%  '
Set-TextValue 'E28' '  -6.72%  '
Set-TextValue 'D29' '8.16'
Set-TextValue 'E29' '  -9.73%  '
Set-TextValue 'E30' '  -3.87%  '
Set-TextValue 'D31' '1.65'
Set-TextValue 'E31' '  -2.56%  '
Set-TextValue 'E32' '  -0.01%  '
Set-TextValue 'D33' '25.65'
Set-TextValue 'E33' '  -4.23%  '
Set-TextValue 'D34' '0.158'
Set-TextValue 'E34' '  -3.04%  '
Set-TextValue 'D35' '6.12'
Set-TextValue 'E35' '  -4.95%  '
Set-TextValue 'B36' 'RenzoRestakedETH'
Set-TextValue 'C36' 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextValue 'D36' '3.516.94'
Set-TextValue 'E36' '  -4.15%  '
Set-TextValue 'B37' 'ImmutableX'
Set-TextValue 'C37' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D37' '1.86'
Set-TextValue 'E37' '  -7.10%  '
Set-TextValue 'D38' '7.99'
Set-TextValue 'E38' '  -5.41%  '
Set-TextValue 'E40' '  +0.09%  '
Set-TextValue 'D41' '174.12'
Set-TextValue 'E41' '  -2.38%  '
Set-TextValue 'B42' 'Filecoin'
Set-TextValue 'C42' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D42' '5.54'
Set-TextValue 'E42' '  -5.65%  '
Set-TextValue 'B43' 'Stacks'
Set-TextValue 'C43' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D43' '2.12'
Set-TextValue 'E43' '  -3.44%  '
Set-TextValue 'D44' '0.0857'
Set-TextValue 'E44' '  -4.55%  '
Set-TextValue 'D45' '0.887'
Set-TextValue 'E45' '  -4.14%  '
Set-TextValue 'D46' '45.19'
Set-TextValue 'E46' '  -4.27%  '
Set-TextValue 'D47' '26.98'
Set-TextValue 'E47' '  -6.77%  '
Set-TextValue 'D48' '2.54'
Set-TextValue 'E48' '  -5.70%  '
Set-TextValue 'E49' '  -0.79%  '
Set-TextValue 'D50' '7.54'
Set-TextValue 'E50' '  -3.45%  '
Set-TextValue 'D51' '1.01'
Set-TextValue 'E51' '  -4.87%  '
